# Adding colors to first row in excel output
# Also introduces a new "Department" column (inserted before the old
# column E), populated with "415 - Phytopathology Lab" for the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column E (shifts E.. right by one) ---
$ws.Columns("E:E").Insert()

# --- New "Department" column header + values ---
$ws.Range("E2").Value = "Department"
for ($r = 3; $r -le 14; $r++) {
    $ws.Cells.Item($r, 5).Value = "415 - Phytopathology Lab"
}

# --- Color the grouped headers in row 1 (by budget/project grouping) ---
$ws.Range("L1:Q1").Interior.Color = 13223074    # 00A2C4C9 - B10/20 AgPlenus
$ws.Range("R1:V1").Interior.Color = 13421812    # 00F4CCCC - B20/20 Lavie-Bio (1)
$ws.Range("W1:Z1").Interior.Color = 13882323    # 00D3D3D3 - B20/20 Lavie-Bio (2)
$ws.Range("AA1:AD1").Interior.Color = 10275833  # 00F9CB9C - B40/20 CPB
$ws.Range("AE1:AJ1").Interior.Color = 13888217  # 00D9EAD3 - B70/20 Biomica + B74/20 Canonic
$ws.Range("AK1:AO1").Interior.Color = 13431551  # 00FFF2CC - B74/20 Canonic (cont.)
$ws.Range("AP1:AY1").Interior.Color = 15254943  # 009FC5E8 - B80/20 PRoduct
$ws.Range("AZ1:BE1").Interior.Color = 8242323   # 0093C47D - B72/20 Casterra
$ws.Range("BF1:BG1").Interior.Color = 14471658  # 00EAD1DC - trailing empty cells

Write-Host "done"
